$d = $word.ActiveDocument

# --- First paragraph: **ID__...__ID** line ---
$para1 = $d.Paragraphs(1)
$pf1 = $para1.Range.ParagraphFormat

# Add a paragraph border (top/left/bottom/right) with 5pt space-from-text,
# matching <w:pBdr><w:top w:space="5"/>...</w:pBdr>
$pf1.Borders.DistanceFromTop = 5
$pf1.Borders.DistanceFromLeft = 5
$pf1.Borders.DistanceFromBottom = 5
$pf1.Borders.DistanceFromRight = 5

# Change left indent from 120 twips (6pt) to 225 twips (11.25pt)
$pf1.LeftIndent = 11.25

# Replace the placeholder id text in the first run only, then remove the
# trailing " " run entirely (it was a second run containing a single space).
$oldId = "**ID__AFFARS_pgi_5301_topic_38__ID**"
$newId = "**ID__AFFARS_AFRC_PGI_5301__ID**"
$full1 = $para1.Range
$idStart = $full1.Start
$idEnd = $idStart + $oldId.Length
$idRange = $d.Range($idStart, $idEnd)
$idRange.Text = $newId

# After the text swap, drop the remaining trailing space run so the
# paragraph contains a single run.
$para1b = $d.Paragraphs(1)
$fullAfter = $para1b.Range
$newIdEnd = $fullAfter.Start + $newId.Length
$trailing = $d.Range($newIdEnd, $fullAfter.End)
$trailing.Text = ""
